# "new file with more columns and data."
#
# The original Sheet1 only had data in columns A:C (row 2/3). This change
# adds three more data cells to row 2 (columns D, E, F -- all literal 1s),
# moves the active selection to E27, and nudges the window/recalc metadata
# to match a freshly-edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data in row 2 -- columns D, E, F.
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

# The sheet's active cell / selection moved from F1 to E27.
$ws.Range("E27").Select()

# Make sure everything (including the new cells) is recalculated.
$excel.CalculateFullRebuild()

# The saved window position in the file moved as part of this edit.
$excel.ActiveWindow.Left = -20
$excel.ActiveWindow.Top = -20
